$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.434937333333333
$ws.Range("H2").Value = 4.304812
$ws.Range("I2").Value = 0.5010808920723563
$ws.Range("J2").Value = 0.5010808920723562
$ws.Range("M2").Value = 5.740110333333334
$ws.Range("N2").Value = 17.220331
$ws.Range("O2").Value = 0.2861925343043439
$ws.Range("P2").Value = 0.2861925343043439
$ws.Range("Q2").Value = 8.236698614752445
$ws.Range("R2").Value = 74.130287532772
$ws.Range("S2").Value = 0.1434056103936691
$ws.Range("T2").Value = 0.143405610393669
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.434937333333333
$ws.Range("H3").Value = 4.304812
$ws.Range("I3").Value = 0.5010808920723563
$ws.Range("J3").Value = 0.5010808920723562
$ws.Range("O3").Value = 0.2917347240316885
$ws.Range("P3").Value = 0.2917347240316885
$ws.Range("Q3").Value = 8.396204335475998
$ws.Range("R3").Value = 75.565839019284
$ws.Range("S3").Value = 0.1461826957662812
$ws.Range("T3").Value = 0.1461826957662812
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.434937333333333
$ws.Range("H4").Value = 4.304812
$ws.Range("I4").Value = 0.5010808920723563
$ws.Range("J4").Value = 0.5010808920723562
$ws.Range("M4").Value = 6.759986
$ws.Range("N4").Value = 20.279958
$ws.Range("O4").Value = 0.3370418707750538
$ws.Range("P4").Value = 0.3370418707750538
$ws.Range("Q4").Value = 9.700156284210667
$ws.Range("R4").Value = 87.301406557896
$ws.Range("S4").Value = 0.1688852412736998
$ws.Range("T4").Value = 0.1688852412736997
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.434937333333333
$ws.Range("H5").Value = 4.304812
$ws.Range("I5").Value = 0.5010808920723563
$ws.Range("J5").Value = 0.5010808920723562
$ws.Range("M5").Value = 1.705448333333333
$ws.Range("N5").Value = 5.116345
$ws.Range("O5").Value = 0.0850308708889137
$ws.Range("P5").Value = 0.0850308708889137
$ws.Range("Q5").Value = 2.447211483571111
$ws.Range("R5").Value = 22.02490335214
$ws.Range("S5").Value = 0.04260734463870623
$ws.Range("T5").Value = 0.04260734463870622
$ws.Range("G6").Value = 0.9964423333333334
$ws.Range("I6").Value = 0.3479582011609289
$ws.Range("J6").Value = 0.3479582011609288
$ws.Range("M6").Value = 5.740110333333334
$ws.Range("N6").Value = 17.220331
$ws.Range("O6").Value = 0.2861925343043439
$ws.Range("P6").Value = 0.2861925343043439
$ws.Range("Q6").Value = 5.719688934137445
$ws.Range("R6").Value = 51.47720040723701
$ws.Range("S6").Value = 0.09958303942222696
$ws.Range("T6").Value = 0.0995830394222269
$ws.Range("G7").Value = 0.9964423333333334
$ws.Range("I7").Value = 0.3479582011609289
$ws.Range("J7").Value = 0.3479582011609288
$ws.Range("O7").Value = 0.2917347240316885
$ws.Range("P7").Value = 0.2917347240316885
$ws.Range("Q7").Value = 5.830452135321
$ws.Range("S7").Value = 0.1015114897902464
$ws.Range("T7").Value = 0.1015114897902463
$ws.Range("G8").Value = 0.9964423333333334
$ws.Range("I8").Value = 0.3479582011609289
$ws.Range("J8").Value = 0.3479582011609288
$ws.Range("M8").Value = 6.759986
$ws.Range("N8").Value = 20.279958
$ws.Range("O8").Value = 0.3370418707750538
$ws.Range("P8").Value = 0.3370418707750538
$ws.Range("Q8").Value = 6.735936223140667
$ws.Range("R8").Value = 60.62342600826601
$ws.Range("S8").Value = 0.117276483070802
$ws.Range("T8").Value = 0.1172764830708019
$ws.Range("G9").Value = 0.9964423333333334
$ws.Range("I9").Value = 0.3479582011609289
$ws.Range("J9").Value = 0.3479582011609288
$ws.Range("M9").Value = 1.705448333333333
$ws.Range("N9").Value = 5.116345
$ws.Range("O9").Value = 0.0850308708889137
$ws.Range("P9").Value = 0.0850308708889137
$ws.Range("Q9").Value = 1.699380916646111
$ws.Range("R9").Value = 15.294428249815
$ws.Range("S9").Value = 0.02958718887765361
$ws.Range("T9").Value = 0.0295871888776536
$ws.Range("G10").Value = 0.4323043333333333
$ws.Range("H10").Value = 1.296913
$ws.Range("I10").Value = 0.150960906766715
$ws.Range("J10").Value = 0.1509609067667149
$ws.Range("M10").Value = 5.740110333333334
$ws.Range("N10").Value = 17.220331
$ws.Range("O10").Value = 0.2861925343043439
$ws.Range("P10").Value = 0.2861925343043439
$ws.Range("Q10").Value = 2.481474570911445
$ws.Range("R10").Value = 22.333271138203
$ws.Range("S10").Value = 0.04320388448844794
$ws.Range("T10").Value = 0.04320388448844792
$ws.Range("G11").Value = 0.4323043333333333
$ws.Range("H11").Value = 1.296913
$ws.Range("I11").Value = 0.150960906766715
$ws.Range("J11").Value = 0.1509609067667149
$ws.Range("O11").Value = 0.2917347240316885
$ws.Range("P11").Value = 0.2917347240316885
$ws.Range("Q11").Value = 2.529528944199
$ws.Range("R11").Value = 22.765760497791
$ws.Range("S11").Value = 0.04404053847516105
$ws.Range("T11").Value = 0.04404053847516104
$ws.Range("G12").Value = 0.4323043333333333
$ws.Range("H12").Value = 1.296913
$ws.Range("I12").Value = 0.150960906766715
$ws.Range("J12").Value = 0.1509609067667149
$ws.Range("M12").Value = 6.759986
$ws.Range("N12").Value = 20.279958
$ws.Range("O12").Value = 0.3370418707750538
$ws.Range("P12").Value = 0.3370418707750538
$ws.Range("Q12").Value = 2.922371241072667
$ws.Range("R12").Value = 26.301341169654
$ws.Range("S12").Value = 0.0508801464305521
$ws.Range("T12").Value = 0.05088014643055208
$ws.Range("G13").Value = 0.4323043333333333
$ws.Range("H13").Value = 1.296913
$ws.Range("I13").Value = 0.150960906766715
$ws.Range("J13").Value = 0.1509609067667149
$ws.Range("M13").Value = 1.705448333333333
$ws.Range("N13").Value = 5.116345
$ws.Range("O13").Value = 0.0850308708889137
$ws.Range("P13").Value = 0.0850308708889137
$ws.Range("Q13").Value = 0.7372727047761111
$ws.Range("R13").Value = 6.635454342985
$ws.Range("S13").Value = 0.01283633737255388
$ws.Range("T13").Value = 0.01283633737255387
